$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "network" column header becomes "method" (performance-bonus payout
# method), row 2's data is unchanged.
$ws.Range("D1").Value = "method"

# Move the active selection to D2, matching the sheet's saved selection.
$ws.Range("D2").Select()
